# Fix <Fichier comptable> generation issue
# Add the missing totals row (row 19) to the "Etat Virement" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A through H hold a single blank-space placeholder on the totals row.
$ws.Range("A19:H19").Value = " "

# Totals for MT brut / Taxe / MT Net columns.
$ws.Range("I19").Value = 272000
$ws.Range("J19").Value = 27900
$ws.Range("K19").Value = 244100

Write-Host "Added totals row 19"
